# Commit: "Add kemerovo, del 8 sites"
# This sheet drops the avto-trend21.ru (F/G) and auto-shop-21.ru (H/I) price
# columns and keeps only alyans-auto.ru (moved from J/K into F/G).
# A couple of data points also change (row 5 and row 15 URLs/price).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a couple of data values that changed along with the column move ---

# Row 5: min_price_url changes to a different listing on alyans-auto.ru
$ws.Range("E5").Value = "https://alyans-auto.ru/auto/auto_18452.html"
$ws.Range("K5").Value = "https://alyans-auto.ru/auto/auto_18452.html"

# Row 15: min_price and min_price_url change to a different listing
$ws.Range("D15").Value = 3214190
$ws.Range("E15").Value = "https://alyans-auto.ru/auto/auto_19834.html"
$ws.Range("J15").Value = 3214190
$ws.Range("K15").Value = "https://alyans-auto.ru/auto/auto_19834.html"

# --- Move the alyans-auto.ru price/url columns (currently J/K) into F/G ---
# (F/G previously held avto-trend21.ru data, which is being removed)

$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 10).Value2   # F = J (price)
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 11).Value2   # G = K (url)
}

# --- Update header row for the surviving F/G columns ---
$ws.Range("F1").Value = "alyans-auto.ru_price"
$ws.Range("G1").Value = "alyans-auto.ru"

# --- Remove the now-obsolete columns: old F/G (avto-trend21.ru), H/I
#     (auto-shop-21.ru) and the now-duplicated J/K (alyans-auto.ru) ---
$ws.Range("H1:K1").EntireColumn.Delete()
